$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 15.89577633333333
$ws.Range("H2").Value = 47.687329
$ws.Range("I2").Value = 0.286059172443548
$ws.Range("J2").Value = 0.2860591724435479
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.477432
$ws.Range("N2").Value = 4.432296
$ws.Range("O2").Value = 0.866150398379125
$ws.Range("P2").Value = 0.866150398379125
$ws.Range("Q2").Value = 23.48492861970933
$ws.Range("R2").Value = 211.364357577384
$ws.Range("S2").Value = 0.2477702661719819
$ws.Range("T2").Value = 0.2477702661719819

# Row 3
$ws.Range("G3").Value = 15.89577633333333
$ws.Range("H3").Value = 47.687329
$ws.Range("I3").Value = 0.286059172443548
$ws.Range("J3").Value = 0.2860591724435479
$ws.Range("O3").Value = 0.133849601620875
$ws.Range("P3").Value = 0.133849601620875
$ws.Range("Q3").Value = 3.629217680584444
$ws.Range("R3").Value = 32.66295912526
$ws.Range("S3").Value = 0.03828890627156609
$ws.Range("T3").Value = 0.03828890627156608

# Row 4
$ws.Range("I4").Value = 0.6735478078679881
$ws.Range("J4").Value = 0.673547807867988
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.477432
$ws.Range("N4").Value = 4.432296
$ws.Range("O4").Value = 0.866150398379125
$ws.Range("P4").Value = 0.866150398379125
$ws.Range("Q4").Value = 55.29702842464533
$ws.Range("R4").Value = 497.673255821808
$ws.Range("S4").Value = 0.5833937021122442
$ws.Range("T4").Value = 0.5833937021122441

# Row 5
$ws.Range("I5").Value = 0.6735478078679881
$ws.Range("J5").Value = 0.673547807867988
$ws.Range("O5").Value = 0.133849601620875
$ws.Range("P5").Value = 0.133849601620875
$ws.Range("R5").Value = 76.90739062612001
$ws.Range("S5").Value = 0.09015410575574385
$ws.Range("T5").Value = 0.09015410575574384

# Row 6
$ws.Range("I6").Value = 0.04039301968846393
$ws.Range("J6").Value = 0.04039301968846393
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.477432
$ws.Range("N6").Value = 4.432296
$ws.Range("O6").Value = 0.866150398379125
$ws.Range("P6").Value = 0.866150398379125
$ws.Range("Q6").Value = 3.316192157080001
$ws.Range("R6").Value = 29.84572941372
$ws.Range("S6").Value = 0.03498643009489887
$ws.Range("T6").Value = 0.03498643009489887

# Row 7
$ws.Range("I7").Value = 0.04039301968846393
$ws.Range("J7").Value = 0.04039301968846393
$ws.Range("O7").Value = 0.133849601620875
$ws.Range("P7").Value = 0.133849601620875
$ws.Range("S7").Value = 0.005406589593565058
$ws.Range("T7").Value = 0.005406589593565057
